# Updated symbol list (Price / Volume(1h) columns) to match the
# Wed Feb 15 23:29:11 UTC 2023 GitHub Actions refresh.
#
# D/E columns are stored as literal text (not numbers/percentages) in the
# source workbook, so each new value is entered with a leading apostrophe
# to force text entry (avoids Excel auto-converting "315.52" -> a number or
# "6.46%" -> a percentage). The Style reset afterwards clears the
# "quote prefix" formatting flag that the apostrophe entry leaves behind,
# so the cell keeps the workbook's default (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'6.46%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'45.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.50%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.175"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.86%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08079"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.50%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.534"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.26%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.680"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.14%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.092"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'17.16%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'8.09%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1929"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.83%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09460"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'5.22%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04297"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'7.43%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.1046"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.52%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001313"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.98%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005974"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'3.21%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D17").Value = "'3.400"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'2.411"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.10%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3368"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.48%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.285"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.86%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-2.37%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.3145"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'4.83%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04260"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5.08%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001282"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.26%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004226"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'6.14%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'9.30%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'11.55%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05456"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'4.69%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.005779"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-4.53%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007728"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.74%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1422"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.57%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007354"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.44%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008577"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'18.46%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3139"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'5.71%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006808"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.30%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-0.41%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.06946"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'51.54%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003985"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-5.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002093"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.41%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001993"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.41%"
$ws.Range("E51").Style = "Normal"
